$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 191 - this pushes the existing rows
# 191:222 down to 193:224 (and extends the used range to A1:R224),
# mirroring Excel's own Rows.Insert() shift-down behaviour.
$ws.Rows("191:192").Insert()

# New weekly entry: Femacal de La Calera / Coquimbo / Zanahoria, split
# into "Primera" and "Segunda" quality rows, dated 44505 (2021-11-05).

# Row 191 - "Primera"
$ws.Cells.Item(191, 1).Value = 3
$ws.Cells.Item(191, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(191, 3).Value = "Coquimbo"
$ws.Cells.Item(191, 4).Value = 44505
$ws.Cells.Item(191, 5).Value = 5
$ws.Cells.Item(191, 6).Value = 100114013
$ws.Cells.Item(191, 7).Value = "Zanahoria"
$ws.Cells.Item(191, 8).Value = "Sin especificar"
$ws.Cells.Item(191, 9).Value = "Primera"
$ws.Cells.Item(191, 10).Value = 230
$ws.Cells.Item(191, 11).Value = 7000
$ws.Cells.Item(191, 12).Value = 7000
$ws.Cells.Item(191, 13).Value = 7000
$ws.Cells.Item(191, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(191, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(191, 16).Value = 350
$ws.Cells.Item(191, 17).Value = 20
$ws.Cells.Item(191, 18).Value = "Hortaliza"

# Row 192 - "Segunda"
$ws.Cells.Item(192, 1).Value = 3
$ws.Cells.Item(192, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(192, 3).Value = "Coquimbo"
$ws.Cells.Item(192, 4).Value = 44505
$ws.Cells.Item(192, 5).Value = 5
$ws.Cells.Item(192, 6).Value = 100114013
$ws.Cells.Item(192, 7).Value = "Zanahoria"
$ws.Cells.Item(192, 8).Value = "Sin especificar"
$ws.Cells.Item(192, 9).Value = "Segunda"
$ws.Cells.Item(192, 10).Value = 250
$ws.Cells.Item(192, 11).Value = 5000
$ws.Cells.Item(192, 12).Value = 5000
$ws.Cells.Item(192, 13).Value = 5000
$ws.Cells.Item(192, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(192, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(192, 16).Value = 250
$ws.Cells.Item(192, 17).Value = 20
$ws.Cells.Item(192, 18).Value = "Hortaliza"
